$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -18.47473578769743
$ws.Cells.Item(2, 3).Value = 2.452238424989163
$ws.Cells.Item(2, 4).Value = -18.47473578769743
$ws.Cells.Item(2, 5).Value = -18.47473578769743
$ws.Cells.Item(2, 6).Value = -18.47473578769743
$ws.Cells.Item(2, 7).Value = -18.47473578769743
$ws.Cells.Item(2, 8).Value = -18.47473578769743
$ws.Cells.Item(2, 9).Value = -18.47473578769743
$ws.Cells.Item(2, 10).Value = -18.47473578769743
$ws.Cells.Item(2, 11).Value = -18.47473578769743

$ws.Cells.Item(3, 2).Value = -18.47473578769743
$ws.Cells.Item(3, 3).Value = -18.47473578769743
$ws.Cells.Item(3, 4).Value = -18.47473578769743
$ws.Cells.Item(3, 5).Value = -18.47473578769743
$ws.Cells.Item(3, 6).Value = -18.47473578769743
$ws.Cells.Item(3, 7).Value = -18.47473578769743
$ws.Cells.Item(3, 8).Value = -18.47473578769743
$ws.Cells.Item(3, 9).Value = 2.303130315053395
$ws.Cells.Item(3, 10).Value = -18.47473578769743
$ws.Cells.Item(3, 11).Value = -18.47473578769743

$ws.Cells.Item(4, 2).Value = -18.47473578769743
$ws.Cells.Item(4, 3).Value = 2.158391596965992
$ws.Cells.Item(4, 4).Value = 2.866305728403292
$ws.Cells.Item(4, 5).Value = -18.47473578769743
$ws.Cells.Item(4, 6).Value = 2.579877141175199
$ws.Cells.Item(4, 7).Value = -18.47473578769743
$ws.Cells.Item(4, 8).Value = 1.863174612261301
$ws.Cells.Item(4, 9).Value = -18.47473578769743
$ws.Cells.Item(4, 10).Value = 2.288564935516257
$ws.Cells.Item(4, 11).Value = -18.47473578769743

$ws.Cells.Item(5, 2).Value = -18.47473578769743
$ws.Cells.Item(5, 3).Value = 0.9448150364698127
$ws.Cells.Item(5, 4).Value = -18.47473578769743
$ws.Cells.Item(5, 5).Value = -18.47473578769743
$ws.Cells.Item(5, 6).Value = -18.47473578769743
$ws.Cells.Item(5, 7).Value = 1.8507113173064
$ws.Cells.Item(5, 8).Value = -18.47473578769743
$ws.Cells.Item(5, 9).Value = -18.47473578769743
$ws.Cells.Item(5, 10).Value = -18.47473578769743
$ws.Cells.Item(5, 11).Value = -18.47473578769743

$ws.Cells.Item(6, 2).Value = -18.47473578769743
$ws.Cells.Item(6, 3).Value = -18.47473578769743
$ws.Cells.Item(6, 4).Value = -18.47473578769743
$ws.Cells.Item(6, 5).Value = -18.47473578769743
$ws.Cells.Item(6, 6).Value = -18.47473578769743
$ws.Cells.Item(6, 7).Value = -18.47473578769743
$ws.Cells.Item(6, 8).Value = -18.47473578769743
$ws.Cells.Item(6, 9).Value = -18.47473578769743
$ws.Cells.Item(6, 10).Value = -18.47473578769743
$ws.Cells.Item(6, 11).Value = -18.47473578769743

$ws.Cells.Item(7, 2).Value = 2.964034241257319
$ws.Cells.Item(7, 3).Value = -18.47473578769743
$ws.Cells.Item(7, 4).Value = -18.47473578769743
$ws.Cells.Item(7, 5).Value = -18.47473578769743
$ws.Cells.Item(7, 6).Value = -18.47473578769743
$ws.Cells.Item(7, 7).Value = -18.47473578769743
$ws.Cells.Item(7, 8).Value = -18.47473578769743
$ws.Cells.Item(7, 9).Value = -18.47473578769743
$ws.Cells.Item(7, 10).Value = -18.47473578769743
$ws.Cells.Item(7, 11).Value = -18.47473578769743

$ws.Cells.Item(8, 2).Value = -18.47473578769743
$ws.Cells.Item(8, 3).Value = -18.47473578769743
$ws.Cells.Item(8, 4).Value = -18.47473578769743
$ws.Cells.Item(8, 5).Value = 2.848150138903017
$ws.Cells.Item(8, 6).Value = -18.47473578769743
$ws.Cells.Item(8, 7).Value = -18.47473578769743
$ws.Cells.Item(8, 8).Value = -18.47473578769743
$ws.Cells.Item(8, 9).Value = -18.47473578769743
$ws.Cells.Item(8, 10).Value = -18.47473578769743
$ws.Cells.Item(8, 11).Value = -18.47473578769743

$ws.Cells.Item(9, 2).Value = 3.608445176423999
$ws.Cells.Item(9, 3).Value = -18.47473578769743
$ws.Cells.Item(9, 4).Value = -18.47473578769743
$ws.Cells.Item(9, 5).Value = -18.47473578769743
$ws.Cells.Item(9, 6).Value = -18.47473578769743
$ws.Cells.Item(9, 7).Value = -18.47473578769743
$ws.Cells.Item(9, 8).Value = -18.47473578769743
$ws.Cells.Item(9, 9).Value = -18.47473578769743
$ws.Cells.Item(9, 10).Value = -18.47473578769743
$ws.Cells.Item(9, 11).Value = -18.47473578769743

$ws.Cells.Item(10, 2).Value = -18.47473578769743
$ws.Cells.Item(10, 3).Value = -18.47473578769743
$ws.Cells.Item(10, 4).Value = -18.47473578769743
$ws.Cells.Item(10, 5).Value = -18.47473578769743
$ws.Cells.Item(10, 6).Value = -18.47473578769743
$ws.Cells.Item(10, 7).Value = -18.47473578769743
$ws.Cells.Item(10, 8).Value = -18.47473578769743
$ws.Cells.Item(10, 9).Value = 1.466095382100302
$ws.Cells.Item(10, 10).Value = -18.47473578769743
$ws.Cells.Item(10, 11).Value = -18.47473578769743

$ws.Cells.Item(11, 2).Value = -18.47473578769743
$ws.Cells.Item(11, 3).Value = -18.47473578769743
$ws.Cells.Item(11, 4).Value = -18.47473578769743
$ws.Cells.Item(11, 5).Value = 1.942610965300625
$ws.Cells.Item(11, 6).Value = -18.47473578769743
$ws.Cells.Item(11, 7).Value = 2.358034497436388
$ws.Cells.Item(11, 8).Value = -18.47473578769743
$ws.Cells.Item(11, 9).Value = -18.47473578769743
$ws.Cells.Item(11, 10).Value = -18.47473578769743
$ws.Cells.Item(11, 11).Value = 4.321924332625012

$ws.Cells.Item(12, 2).Value = -18.47473578769743
$ws.Cells.Item(12, 3).Value = -18.47473578769743
$ws.Cells.Item(12, 4).Value = -18.47473578769743
$ws.Cells.Item(12, 5).Value = -18.47473578769743
$ws.Cells.Item(12, 6).Value = -18.47473578769743
$ws.Cells.Item(12, 7).Value = -18.47473578769743
$ws.Cells.Item(12, 8).Value = -18.47473578769743
$ws.Cells.Item(12, 9).Value = -18.47473578769743
$ws.Cells.Item(12, 10).Value = -18.47473578769743
$ws.Cells.Item(12, 11).Value = -18.47473578769743

$ws.Cells.Item(13, 2).Value = -18.47473578769743
$ws.Cells.Item(13, 3).Value = -18.47473578769743
$ws.Cells.Item(13, 4).Value = -18.47473578769743
$ws.Cells.Item(13, 5).Value = 1.660974750745874
$ws.Cells.Item(13, 6).Value = -18.47473578769743
$ws.Cells.Item(13, 7).Value = -18.47473578769743
$ws.Cells.Item(13, 8).Value = -18.47473578769743
$ws.Cells.Item(13, 9).Value = -18.47473578769743
$ws.Cells.Item(13, 10).Value = 2.340172851074819
$ws.Cells.Item(13, 11).Value = -18.47473578769743

$ws.Cells.Item(14, 2).Value = -18.47473578769743
$ws.Cells.Item(14, 3).Value = -18.47473578769743
$ws.Cells.Item(14, 4).Value = 1.643309251894609
$ws.Cells.Item(14, 5).Value = -18.47473578769743
$ws.Cells.Item(14, 6).Value = -18.47473578769743
$ws.Cells.Item(14, 7).Value = -18.47473578769743
$ws.Cells.Item(14, 8).Value = -18.47473578769743
$ws.Cells.Item(14, 9).Value = -18.47473578769743
$ws.Cells.Item(14, 10).Value = -18.47473578769743
$ws.Cells.Item(14, 11).Value = -18.47473578769743

$ws.Cells.Item(15, 2).Value = -18.47473578769743
$ws.Cells.Item(15, 3).Value = -18.47473578769743
$ws.Cells.Item(15, 4).Value = -0.2848485059878001
$ws.Cells.Item(15, 5).Value = -18.47473578769743
$ws.Cells.Item(15, 6).Value = -18.47473578769743
$ws.Cells.Item(15, 7).Value = -18.47473578769743
$ws.Cells.Item(15, 8).Value = -18.47473578769743
$ws.Cells.Item(15, 9).Value = -18.47473578769743
$ws.Cells.Item(15, 10).Value = -18.47473578769743
$ws.Cells.Item(15, 11).Value = -18.47473578769743

$ws.Cells.Item(16, 2).Value = -18.47473578769743
$ws.Cells.Item(16, 3).Value = -18.47473578769743
$ws.Cells.Item(16, 4).Value = -18.47473578769743
$ws.Cells.Item(16, 5).Value = -18.47473578769743
$ws.Cells.Item(16, 6).Value = -18.47473578769743
$ws.Cells.Item(16, 7).Value = -18.47473578769743
$ws.Cells.Item(16, 8).Value = -18.47473578769743
$ws.Cells.Item(16, 9).Value = -18.47473578769743
$ws.Cells.Item(16, 10).Value = 2.471288774060023
$ws.Cells.Item(16, 11).Value = -18.47473578769743

$ws.Cells.Item(17, 2).Value = -18.47473578769743
$ws.Cells.Item(17, 3).Value = 0.627375765368265
$ws.Cells.Item(17, 4).Value = -0.1189502618190623
$ws.Cells.Item(17, 5).Value = -18.47473578769743
$ws.Cells.Item(17, 6).Value = -18.47473578769743
$ws.Cells.Item(17, 7).Value = -18.47473578769743
$ws.Cells.Item(17, 8).Value = 0.2440679696031669
$ws.Cells.Item(17, 9).Value = 0.6115481873847308
$ws.Cells.Item(17, 10).Value = 1.190692085595782
$ws.Cells.Item(17, 11).Value = -18.47473578769743

$ws.Cells.Item(18, 2).Value = -18.47473578769743
$ws.Cells.Item(18, 3).Value = -18.47473578769743
$ws.Cells.Item(18, 4).Value = -18.47473578769743
$ws.Cells.Item(18, 5).Value = -18.47473578769743
$ws.Cells.Item(18, 6).Value = -18.47473578769743
$ws.Cells.Item(18, 7).Value = -18.47473578769743
$ws.Cells.Item(18, 8).Value = 0.4253910201745995
$ws.Cells.Item(18, 9).Value = 0.8691720836451482
$ws.Cells.Item(18, 10).Value = 1.152258310334972
$ws.Cells.Item(18, 11).Value = -18.47473578769743

$ws.Cells.Item(19, 2).Value = -18.47473578769743
$ws.Cells.Item(19, 3).Value = -18.47473578769743
$ws.Cells.Item(19, 4).Value = 1.684082197847617
$ws.Cells.Item(19, 5).Value = -18.47473578769743
$ws.Cells.Item(19, 6).Value = -18.47473578769743
$ws.Cells.Item(19, 7).Value = -18.47473578769743
$ws.Cells.Item(19, 8).Value = 1.888339144406526
$ws.Cells.Item(19, 9).Value = 2.118583293380931
$ws.Cells.Item(19, 10).Value = -18.47473578769743
$ws.Cells.Item(19, 11).Value = -18.47473578769743

$ws.Cells.Item(20, 2).Value = -18.47473578769743
$ws.Cells.Item(20, 3).Value = 1.706073005745431
$ws.Cells.Item(20, 4).Value = 2.210746189294917
$ws.Cells.Item(20, 5).Value = -18.47473578769743
$ws.Cells.Item(20, 6).Value = 3.809523802324861
$ws.Cells.Item(20, 7).Value = -18.47473578769743
$ws.Cells.Item(20, 8).Value = 2.319994980152691
$ws.Cells.Item(20, 9).Value = 2.203129207449515
$ws.Cells.Item(20, 10).Value = -18.47473578769743
$ws.Cells.Item(20, 11).Value = -18.47473578769743

$ws.Cells.Item(21, 2).Value = -18.47473578769743
$ws.Cells.Item(21, 3).Value = 1.735841627572201
$ws.Cells.Item(21, 4).Value = -18.47473578769743
$ws.Cells.Item(21, 5).Value = 2.534281572140444
$ws.Cells.Item(21, 6).Value = -18.47473578769743
$ws.Cells.Item(21, 7).Value = 3.49396229021021
$ws.Cells.Item(21, 8).Value = 2.361552701152363
$ws.Cells.Item(21, 9).Value = -18.47473578769743
$ws.Cells.Item(21, 10).Value = -18.47473578769743
$ws.Cells.Item(21, 11).Value = -18.47473578769743
